$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = 173830
$ws.Range("C4").Value = 163899
$ws.Range("C5").Value = 9931
$ws.Range("C6").Value = 771
$ws.Range("C7").Value = 5.71
$ws.Range("C8").Value = 64.35
